$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of accelerometer data to insert right after the header row (row 1),
# pushing the existing data rows down by 8.
$dataTop = @(
    @(-1.275631546974182, 1.748281717300415, 0.7527783811092381),
    @(-1.113769233226775, 1.769958406686783, 1.082688376307487),
    @(0.6056947708129873, 1.413846492767334, 1.041245818138122),
    @(-0.2511940002441411, 1.83873063325882, 0.5010688602924345),
    @(-0.7442607879638676, 1.761505782604217, 0.9367214739322666),
    @(-0.8476336002349854, 1.69824892282486, 0.9451412782073016),
    @(-0.6413483619689934, 1.650843858718872, 0.9322790801525122),
    @(-0.5721501111984255, 1.609763711690903, 1.015784159302711)
)

# New rows of accelerometer data to append after the (now shifted) last row.
$dataBottom = @(
    @(-0.4257340431213359, 1.845006287097933, 0.8945446908474008),
    @(0.06665813922882169, 1.836877554655074, 0.7217497229576104)
)

$insertCount = $dataTop.Count

# Insert blank rows right below the header (row 1) and strip any formatting
# that Excel copies in from neighboring rows so the new cells stay unstyled,
# matching the rest of the data rows.
$ws.Rows.Item(2).Resize($insertCount).Insert()
$ws.Rows.Item(2).Resize($insertCount).ClearFormats()

for ($i = 0; $i -lt $dataTop.Count; $i++) {
    $r = 2 + $i
    $row = $dataTop[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Find the new last used row after the insert, then append the extra rows.
$lastRow = $ws.UsedRange.Rows.Count

for ($i = 0; $i -lt $dataBottom.Count; $i++) {
    $r = $lastRow + 1 + $i
    $row = $dataBottom[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

Write-Output "done"
